# Added DynamicListingPage_Test; removed unneeded config files
$wb = $excel.ActiveWorkbook

# --- Update the BlogPostPage cursor position (cosmetic selection move) ---
$blogPost = $wb.Worksheets.Item("BlogPostPage")
$blogPost.Range("J35").Select() | Out-Null

# --- Insert the new "DynamicListingPage" sheet right after "CTHPPage" ---
$cthp = $wb.Worksheets.Item("CTHPPage")
$ws = $wb.Worksheets.Add($null, $cthp)
$ws.Name = "DynamicListingPage"

# Copy header formatting (bold + fill) from an existing sheet's header row
$headerSrc = $wb.Worksheets.Item("CTHPPage").Range("A1:B1")
$headerSrc.Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# --- Fill in the data ---
# (Values are written in the same cell-by-cell order the sheet was
# originally authored in, so newly introduced shared strings land at the
# same indices as the source workbook: ContentType before Path within
# each new disease/intervention/manual group.)
$ws.Range("A1").Value = "Path"
$ws.Range("B1").Value = "ContentType"

$ws.Range("B2").Value = "Disease Listing Page"
$ws.Range("A2").Value = "/about-cancer/treatment/clinical-trials/disease/breast-cancer"
$ws.Range("A3").Value = "/about-cancer/treatment/clinical-trials/disease/breast-cancer/treatment"
$ws.Range("A4").Value = "/about-cancer/treatment/clinical-trials/disease/breast-cancer/treatment/trastuzumab"
$ws.Range("B3").Value = "Disease Listing Page"
$ws.Range("B4").Value = "Disease Listing Page"

$ws.Range("B5").Value = "Intervention Listing Page"
$ws.Range("A5").Value = "/about-cancer/treatment/clinical-trials/intervention/trastuzumab"
$ws.Range("A6").Value = "/about-cancer/treatment/clinical-trials/intervention/trastuzumab/treatment"
$ws.Range("B6").Value = "Intervention Listing Page"

$ws.Range("B7").Value = "Manual Listing Page"
$ws.Range("A7").Value = "/about-cancer/treatment/clinical-trials/kidney-cancer"

# --- Column widths to match the bestFit sizing used on the other sheets ---
$ws.Columns.Item(1).ColumnWidth = 79.1666666666667
$ws.Columns.Item(2).ColumnWidth = 22.5924479166667

# --- Cursor / active-cell position on the new sheet ---
$ws.Range("A9").Select() | Out-Null
